$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new question/answer rows (14-22) below the existing data (rows 2-13).
# This also grows the sheet dimension from A2:B13 to A2:B22 automatically.

$ws.Range("A14").Value = "Hi"
$ws.Range("B14").Value = "Hello! How can I assist you today?"

$ws.Range("A15").Value = "Can you tell me about the fatigue"
$ws.Range("B15").Value = "Fatigue is a common symptom reported by patients with rheumatoid arthritis (RA). Patients often experience different forms of fatigue, including mental and physical fatigue, fatigue with or without a reason (such as poor sleep or physical activity), fatigue with or without pain, and fatigue with or without dizziness. Some patients also describe fatigue with the desire to go to bed and sleep, while others feel tired but don't have the desire to sleep. Fatigue can vary in severity, frequency, and duration, making it difficult for patients to specify these aspects. `nA definition of fatigue in RA has been provided by Dittner et al., stating that it is an essentially subjective experience characterized by extreme and persistent tiredness, weakness, or exhaustion, both mentally and physically. However, not all patients experience fatigue in the same way. While most patients report a negative impact of fatigue on their lives, the severity of fatigue can vary among individuals.`nFatigue in RA is a complex phenomenon influenced by physical, psychological, and social resources and competencies. It is important to understand the individual experiences and perspectives of patients with fatigue in order to provide appropriate support and management strategies."

$ws.Range("A16").Value = "You bro I am so tired, give me advice"
$ws.Range("B16").Value = "It sounds like you're feeling really tired. Here are a few suggestions that might help:`n1. Make sure you're getting enough sleep. Establish a consistent bedtime routine and create a sleep-friendly environment by turning off electronics before bed.`n2. Take breaks and rest between activities to allow your muscles to refuel. It's important to spread out your activities throughout the day and plan more strenuous tasks for times when you typically have the most energy.`n3. Don't be afraid to ask for help. If you're feeling overwhelmed and don't have the energy to do everything, reach out to family members or consider hiring assistance for tasks like housecleaning.`n4. Consider utilizing pain management tools and resources to help manage any underlying pain that might be contributing to your fatigue.`nRemember, I'm just an AI and these are general suggestions. If your fatigue persists or worsens, it's always a good idea to consult with a healthcare professional for personalized advice."

$ws.Range("A17").Value = "Who is the president of the Kazakstan?"
$ws.Range("B17").Value = "I don't know the answer to that question based on the provided context."

$ws.Range("A18").Value = "Who is the president of US?"
$ws.Range("B18").Value = "I don't have access to real-time information, but as of my last update in 2009, Barack Obama was the President of the United States."

$ws.Range("A19").Value = "I want to eat something, which things can help me to destroy fatigue?"
$ws.Range("B19").Value = "Eating the right foods can help combat fatigue. Some options include fresh fruits and vegetables, lean proteins, dairy products, and healthy fats. Foods rich in omega-3 fatty acids, like nuts and cold-water fish, can also help ease fatigue by reducing inflammation. Eating small, healthy meals throughout the day and choosing snacks like an apple with peanut butter or tuna on whole-grain bread can provide sustained energy. Additionally, starting your day with a breakfast that includes high-fiber cereal, fresh fruit, and low-fat milk, or adding a hard-boiled egg, a cup of yogurt, or a banana to your morning can give you a boost of energy. Remember to stay hydrated by drinking water throughout the day as well."

$ws.Range("A20").Value = "Which year Israel is founded?"
$ws.Range("B20").Value = "The context provided does not mention the year in which Israel was founded. Therefore, I do not have the information to answer your question."

$ws.Range("A21").Value = "Which year Palestine is founded?"
$ws.Range("B21").Value = "Palestine is a region that has a long history and has been inhabited for thousands of years. The modern political entity known as Palestine was founded in 1988, when the Palestinian National Council declared the establishment of the State of Palestine. However, it is important to note that the recognition and political status of Palestine is a complex and ongoing issue."

$ws.Range("A22").Value = "Your bro"
$ws.Range("B22").Value = "I'm sorry, but I don't understand what you mean by `"Your bro.`" Can you please provide more context or clarify your question?"
